$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 764, shifting existing rows 764-832 down to 765-833
$ws.Rows(764).Insert()

# Populate the new row 764 with the new record's data
$ws.Cells.Item(764, 1).Value2  = 3
$ws.Cells.Item(764, 2).Value2  = "Femacal de La Calera"
$ws.Cells.Item(764, 3).Value2  = "Coquimbo"
$ws.Cells.Item(764, 4).Value2  = 45166
$ws.Cells.Item(764, 5).Value2  = 5
$ws.Cells.Item(764, 6).Value2  = 100112032
$ws.Cells.Item(764, 7).Value2  = "Zapallo italiano"
$ws.Cells.Item(764, 8).Value2  = "Sin especificar"
$ws.Cells.Item(764, 9).Value2  = "Primera"
$ws.Cells.Item(764, 10).Value2 = 80
$ws.Cells.Item(764, 11).Value2 = 14000
$ws.Cells.Item(764, 12).Value2 = 15000
$ws.Cells.Item(764, 13).Value2 = 14500
$ws.Cells.Item(764, 14).Value2 = "$/caja 60 unidades"
$ws.Cells.Item(764, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(764, 16).Value2 = 242
$ws.Cells.Item(764, 17).Value2 = 60
$ws.Cells.Item(764, 18).Value2 = "Hortaliza"

# Make sure the date cell keeps the date number-format style used by the rest of column D
$ws.Cells.Item(764, 4).NumberFormat = $ws.Cells.Item(765, 4).NumberFormat
